# Update "想去人数" (number of people interested) counts that changed
# between two scrapes of the 合肥-漫展信息 data.
#
# 展览 (Exhibitions) sheet:
#   F2: 5437 -> 5439
#   F6: 819  -> 820
#   F7: 19   -> 20
#
# 演出 (Performances) sheet:
#   F3: 15 -> 16
#
# 全部类型 (All types) sheet - aggregated view of the above rows:
#   F2:  5437 -> 5439
#   F6:  819  -> 820
#   F7:  19   -> 20
#   F11: 15   -> 16

$wb = $excel.ActiveWorkbook

$wsExhibitions = $wb.Worksheets.Item("展览")
$wsExhibitions.Range("F2").Value = 5439
$wsExhibitions.Range("F6").Value = 820
$wsExhibitions.Range("F7").Value = 20

$wsPerformances = $wb.Worksheets.Item("演出")
$wsPerformances.Range("F3").Value = 16

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 5439
$wsAll.Range("F6").Value = 820
$wsAll.Range("F7").Value = 20
$wsAll.Range("F11").Value = 16
